$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add a new column P (year 2022) mirroring column O's layout/formatting
#    for every existing row (2-6; row 1 keeps its original span only).
# ---------------------------------------------------------------------------
$ws.Range("O2").Copy() | Out-Null
$ws.Range("P2").PasteSpecial(-4122) | Out-Null

$ws.Range("O3").Copy() | Out-Null
$ws.Range("P3").PasteSpecial(-4122) | Out-Null

$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4").PasteSpecial(-4122) | Out-Null

$ws.Range("O5").Copy() | Out-Null
$ws.Range("P5").PasteSpecial(-4122) | Out-Null

$ws.Range("O6").Copy() | Out-Null
$ws.Range("P6").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 2. Populate the new column's data.
# ---------------------------------------------------------------------------
$ws.Range("P3").Value = 2022
$ws.Range("P4").Formula = "=P5/P6*1000"
$ws.Range("P5").Value = 1339.6
$ws.Range("P6").Value = 6300.5

# ---------------------------------------------------------------------------
# 3. Re-format the numeric indicator rows (4-6) with a single-decimal custom
#    number format (0.0) spanning the original data plus the new column P.
# ---------------------------------------------------------------------------
$ws.Range("D4:P4").NumberFormat = "0.0"
$ws.Range("D5:P5").NumberFormat = "0.0"
$ws.Range("D6:P6").NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# 4. Move the active selection like the source workbook does.
# ---------------------------------------------------------------------------
$ws.Range("S4").Select() | Out-Null
